# The deck's slide master ("Design") was re-themed from the green "Integral"
# theme to the default "Office Theme" colour palette.
#
# The 12 DrawingML theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) are reachable - and genuinely persisted back to the underlying
# theme part - through Slide.ThemeColorScheme(i).RGB, so we drive the swap
# through that object model surface rather than touching the OOXML parts
# directly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index -> new RGB long, as produced by VBA's RGB(r,g,b) = r + g*256 + b*65536.
# Target palette = the standard Office Theme colours.
$newColors = @{
    1  = 0        # dk1      #000000
    2  = 16777215 # lt1      #FFFFFF
    3  = 6968388   # dk2      #44546A
    4  = 15132391  # lt2      #E7E6E6
    5  = 13998939  # accent1  #5B9BD5
    6  = 3243501    # accent2  #ED7D31
    7  = 10855845   # accent3  #A5A5A5
    8  = 49407      # accent4  #FFC000
    9  = 12874308   # accent5  #4472C4
    10 = 4697456    # accent6  #70AD47
    11 = 12673797   # hlink    #0563C1
    12 = 7491477    # folHlink #954F72
}

for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $newColors[$i]
}
